$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cost updates
$ws.Range("D11").Value = 5.99
$ws.Range("D12").Value = 4.96

# New row 14: Heat-shrink tubing
$ws.Range("A14").Value = "Heat-shrink tubing, 50mm"
$ws.Range("C14").Value = "Amazon"
$ws.Range("D14").Value = 5.43
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 30

# Bold/highlight the TOTAL row to match the header style
$ws.Range("A17:G17").Font.Bold = $true

# Update the active selection to match the saved view state
$ws.Range("F19").Select()
